$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 3.584715618562996
$ws.Cells.Item(2, 3).Value = 0.7124868897915917
$ws.Cells.Item(2, 4).Value = 0.07670995491601218
$ws.Cells.Item(2, 5).Value = 0.01715278594704372
$ws.Cells.Item(2, 7).Value = 0.002576924206532098
$ws.Cells.Item(2, 9).Value = 2.33792461914139
$ws.Cells.Item(2, 12).Value = 0.3463167753486687
$ws.Cells.Item(2, 14).Value = 2.36834847398238
$ws.Cells.Item(3, 2).Value = 3.409791459777011
$ws.Cells.Item(3, 3).Value = 0.649435907127156
$ws.Cells.Item(3, 4).Value = 0.06981875351138456
$ws.Cells.Item(3, 5).Value = 0.01676052461432231
$ws.Cells.Item(3, 7).Value = 0.002584884137783769
$ws.Cells.Item(3, 9).Value = 2.298994561044381
$ws.Cells.Item(3, 12).Value = 0.3345601425693019
$ws.Cells.Item(3, 14).Value = 2.36650016725136
$ws.Cells.Item(4, 2).Value = 3.304972070762403
$ws.Cells.Item(4, 3).Value = 0.6111925746459406
$ws.Cells.Item(4, 4).Value = 0.06563747383457041
$ws.Cells.Item(4, 5).Value = 0.01651634184707262
$ws.Cells.Item(4, 7).Value = 0.002590019711914158
$ws.Cells.Item(4, 9).Value = 2.276372573284661
$ws.Cells.Item(4, 12).Value = 0.3275767838947417
$ws.Cells.Item(4, 14).Value = 2.366031438250417
$ws.Cells.Item(5, 2).Value = 3.262899589111157
$ws.Cells.Item(5, 3).Value = 0.5957231724231633
$ws.Cells.Item(5, 4).Value = 0.06394579063957906
$ws.Cells.Item(5, 5).Value = 0.01641597210753787
$ws.Cells.Item(5, 7).Value = 0.002592175156307088
$ws.Cells.Item(5, 9).Value = 2.267473042003161
$ws.Cells.Item(5, 12).Value = 0.3247896098815204
$ws.Cells.Item(5, 14).Value = 2.366006437633075
$ws.Cells.Item(6, 2).Value = 3.255952055657701
$ws.Cells.Item(6, 3).Value = 0.5931613478694544
$ws.Cells.Item(6, 4).Value = 0.06366561637783263
$ws.Cells.Item(6, 5).Value = 0.01639925291774436
$ws.Cells.Item(6, 7).Value = 0.002592536857886048
$ws.Cells.Item(6, 9).Value = 2.266014466711965
$ws.Cells.Item(6, 12).Value = 0.3243303250867058
$ws.Cells.Item(6, 14).Value = 2.366012269527459
$ws.Cells.Item(7, 2).Value = 3.304402076060512
$ws.Cells.Item(7, 3).Value = 0.6109834871369912
$ws.Cells.Item(7, 4).Value = 0.06561461015702719
$ws.Cells.Item(7, 5).Value = 0.0165149917543097
$ws.Cells.Item(7, 7).Value = 0.002590048527042921
$ws.Cells.Item(7, 9).Value = 2.276251263020754
$ws.Cells.Item(7, 12).Value = 0.3275389585430304
$ws.Cells.Item(7, 14).Value = 2.366030430872542
$ws.Cells.Item(8, 2).Value = 3.523860568977966
$ws.Cells.Item(8, 3).Value = 0.6906473697472961
$ws.Cells.Item(8, 4).Value = 0.07432329587261677
$ws.Cells.Item(8, 5).Value = 0.01701821056465747
$ws.Cells.Item(8, 7).Value = 0.002579617444389534
$ws.Cells.Item(8, 9).Value = 2.324233742494854
$ws.Cells.Item(8, 12).Value = 0.3422139167862639
$ws.Cells.Item(8, 14).Value = 2.367572012185533
$ws.Cells.Item(9, 2).Value = 3.975105811234641
$ws.Cells.Item(9, 3).Value = 0.8507521972218228
$ws.Cells.Item(9, 4).Value = 0.09181357296449733
$ws.Cells.Item(9, 5).Value = 0.01797976129167811
$ws.Cells.Item(9, 7).Value = 0.002561119217767581
$ws.Cells.Item(9, 9).Value = 2.428640254854301
$ws.Cells.Item(9, 12).Value = 0.3728856257973661
$ws.Cells.Item(9, 14).Value = 2.3759499877564
$ws.Cells.Item(10, 2).Value = 4.319941277953262
$ws.Cells.Item(10, 3).Value = 0.970975693340165
$ws.Cells.Item(10, 4).Value = 0.1049392873770074
$ws.Cells.Item(10, 5).Value = 0.01867251963511674
$ws.Cells.Item(10, 7).Value = 0.002548705051328908
$ws.Cells.Item(10, 9).Value = 2.511851558253156
$ws.Cells.Item(10, 12).Value = 0.3966163120191339
$ws.Cells.Item(10, 14).Value = 2.385468838339165
$ws.Cells.Item(11, 2).Value = 4.47983448044306
$ws.Cells.Item(11, 3).Value = 1.026282647888479
$ws.Cells.Item(11, 4).Value = 0.1109757609534086
$ws.Cells.Item(11, 5).Value = 0.01898507967868746
$ws.Cells.Item(11, 7).Value = 0.002543309409350028
$ws.Cells.Item(11, 9).Value = 2.551167412840883
$ws.Cells.Item(11, 12).Value = 0.4076811398072095
$ws.Cells.Item(11, 14).Value = 2.390551525195491
$ws.Cells.Item(12, 2).Value = 4.5408270640811
$ws.Cells.Item(12, 3).Value = 1.04731863070748
$ws.Cells.Item(12, 4).Value = 0.1132714622368667
$ws.Cells.Item(12, 5).Value = 0.01910309782832975
$ws.Cells.Item(12, 7).Value = 0.00254130212559551
$ws.Cells.Item(12, 9).Value = 2.566269494590074
$ws.Cells.Item(12, 12).Value = 0.4119106046595675
$ws.Cells.Item(12, 14).Value = 2.392586217727086
$ws.Cells.Item(13, 2).Value = 4.527671298492464
$ws.Cells.Item(13, 3).Value = 1.042783983344691
$ws.Cells.Item(13, 4).Value = 0.1127765987219647
$ws.Cells.Item(13, 5).Value = 0.01907769522523495
$ws.Cells.Item(13, 7).Value = 0.00254173283637703
$ws.Cells.Item(13, 9).Value = 2.563007412604208
$ws.Cells.Item(13, 12).Value = 0.4109979478878643
$ws.Cells.Item(13, 14).Value = 2.392143089776766
$ws.Cells.Item(14, 2).Value = 4.484843413918497
$ws.Cells.Item(14, 3).Value = 1.028011416371157
$ws.Cells.Item(14, 4).Value = 0.1111644305643154
$ws.Cells.Item(14, 5).Value = 0.01899479577774166
$ws.Cells.Item(14, 7).Value = 0.002543143550580035
$ws.Cells.Item(14, 9).Value = 2.552405559408896
$ws.Cells.Item(14, 12).Value = 0.40802830571198
$ws.Cells.Item(14, 14).Value = 2.390716705920823
$ws.Cells.Item(15, 2).Value = 4.458668267746475
$ws.Cells.Item(15, 3).Value = 1.018974944027036
$ws.Cells.Item(15, 4).Value = 0.1101782218080558
$ws.Cells.Item(15, 5).Value = 0.01894397387374536
$ws.Cells.Item(15, 7).Value = 0.002544012325052431
$ws.Cells.Item(15, 9).Value = 2.545939601274114
$ws.Cells.Item(15, 12).Value = 0.4062144729097383
$ws.Cells.Item(15, 14).Value = 2.389857381832968
$ws.Cells.Item(16, 2).Value = 4.309553565165743
$ws.Cells.Item(16, 3).Value = 0.9673739825688017
$ws.Cells.Item(16, 4).Value = 0.1045461418852369
$ws.Cells.Item(16, 5).Value = 0.01865204369364015
$ws.Cells.Item(16, 7).Value = 0.00254906271078347
$ws.Cells.Item(16, 9).Value = 2.509311920460803
$ws.Cells.Item(16, 12).Value = 0.39589868296666
$ws.Cells.Item(16, 14).Value = 2.385151966590996
$ws.Cells.Item(17, 2).Value = 4.21885811038976
$ws.Cells.Item(17, 3).Value = 0.935879060238392
$ws.Cells.Item(17, 4).Value = 0.1011081049491338
$ws.Cells.Item(17, 5).Value = 0.01847231405419425
$ws.Cells.Item(17, 7).Value = 0.002552225224930161
$ws.Cells.Item(17, 9).Value = 2.487219240100472
$ws.Cells.Item(17, 12).Value = 0.3896398167284048
$ws.Cells.Item(17, 14).Value = 2.382459338678004
$ws.Cells.Item(18, 2).Value = 4.16697671127281
$ws.Cells.Item(18, 3).Value = 0.9178218288915332
$ws.Cells.Item(18, 4).Value = 0.09913677976683744
$ws.Cells.Item(18, 5).Value = 0.0183686937589993
$ws.Cells.Item(18, 7).Value = 0.002554067918882623
$ws.Cells.Item(18, 9).Value = 2.474649556658221
$ws.Cells.Item(18, 12).Value = 0.3860652182853528
$ws.Cells.Item(18, 14).Value = 2.38098133215982
$ws.Cells.Item(19, 2).Value = 4.149459095374027
$ws.Cells.Item(19, 3).Value = 0.9117177772314449
$ws.Cells.Item(19, 4).Value = 0.09847036597430758
$ws.Cells.Item(19, 5).Value = 0.01833356682145215
$ws.Cells.Item(19, 7).Value = 0.002554695901556181
$ws.Cells.Item(19, 9).Value = 2.470417172280605
$ws.Cells.Item(19, 12).Value = 0.3848592516008722
$ws.Cells.Item(19, 14).Value = 2.380492997405483
$ws.Cells.Item(20, 2).Value = 4.228483312697563
$ws.Cells.Item(20, 3).Value = 0.9392257378867157
$ws.Cells.Item(20, 4).Value = 0.1014734512462212
$ws.Cells.Item(20, 5).Value = 0.01849147168171417
$ws.Cells.Item(20, 7).Value = 0.002551886118912439
$ws.Cells.Item(20, 9).Value = 2.489556792404244
$ws.Cells.Item(20, 12).Value = 0.390303456663105
$ws.Cells.Item(20, 14).Value = 2.382738641299966
$ws.Cells.Item(21, 2).Value = 4.497410862139191
$ws.Cells.Item(21, 3).Value = 1.032347938197177
$ws.Cells.Item(21, 4).Value = 0.1116376937106622
$ws.Cells.Item(21, 5).Value = 0.01901915438766544
$ws.Cells.Item(21, 7).Value = 0.002542728216470746
$ws.Cells.Item(21, 9).Value = 2.555513742311248
$ws.Cells.Item(21, 12).Value = 0.40889948527456
$ws.Cells.Item(21, 14).Value = 2.39113267000954
$ws.Cells.Item(22, 2).Value = 4.675765976915727
$ws.Cells.Item(22, 3).Value = 1.093749295746591
$ws.Cells.Item(22, 4).Value = 0.1183380420049787
$ws.Cells.Item(22, 5).Value = 0.01936204877892678
$ws.Cells.Item(22, 7).Value = 0.002536952298844248
$ws.Cells.Item(22, 9).Value = 2.599869677914029
$ws.Cells.Item(22, 12).Value = 0.4212833678134587
$ws.Cells.Item(22, 14).Value = 2.397260604601513
$ws.Cells.Item(23, 2).Value = 4.580333845923178
$ws.Cells.Item(23, 3).Value = 1.060927503570952
$ws.Cells.Item(23, 4).Value = 0.1147565514445432
$ws.Cells.Item(23, 5).Value = 0.01917921060058347
$ws.Cells.Item(23, 7).Value = 0.002540015947947563
$ws.Cells.Item(23, 9).Value = 2.576080528613403
$ws.Cells.Item(23, 12).Value = 0.4146525531628242
$ws.Cells.Item(23, 14).Value = 2.393930675418289
$ws.Cells.Item(24, 2).Value = 4.224130947648746
$ws.Cells.Item(24, 3).Value = 0.9377125506532593
$ws.Cells.Item(24, 4).Value = 0.1013082618692636
$ws.Cells.Item(24, 5).Value = 0.01848281142785968
$ws.Cells.Item(24, 7).Value = 0.002552039352468722
$ws.Cells.Item(24, 9).Value = 2.488499575092774
$ws.Cells.Item(24, 12).Value = 0.3900033512578744
$ws.Cells.Item(24, 14).Value = 2.382612150618613
$ws.Cells.Item(25, 2).Value = 3.850737555175044
$ws.Cells.Item(25, 3).Value = 0.8070006309040991
$ws.Cells.Item(25, 4).Value = 0.08703535973486964
$ws.Cells.Item(25, 5).Value = 0.01772216629644685
$ws.Cells.Item(25, 7).Value = 0.002565915674529324
$ws.Cells.Item(25, 9).Value = 2.399269320998911
$ws.Cells.Item(25, 12).Value = 0.3643811083563264
$ws.Cells.Item(25, 14).Value = 2.373100110409055
